$d = $word.ActiveDocument

function TrimCR([string]$s) {
    return $s.TrimEnd([char]13)
}

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the Heading1 title
#    paragraph. We clone the existing bold "Play Desperate Dawgs..." paragraph
#    (near the end of the doc) via Copy/Paste so the new paragraph picks up
#    the same run layout (leading empty run + bold run), then edit its text
#    and append the non-bold remainder.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titleText = "Play Desperate Dawgs 2 Gigablox Free - Review 2021"

# Locate the paragraph that currently holds the bold
# "Play Desperate Dawgs 2 Gigablox Free - Review 2021" text near the end of
# the document (it still has the original run structure we want to reuse).
$sourcePara = $null
for ($i = 2; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ((TrimCR $para.Range.Text) -eq $titleText) {
        $sourcePara = $para
    }
}

$sourcePara.Range.Copy()

$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"
$metaPara2 = $d.Paragraphs(2)
$metaPara2.Range.Paste()

$metaPara3 = $d.Paragraphs(2)
$boldRange = $d.Range($metaPara3.Range.Start, $metaPara3.Range.End - 1)
$boldRange.Text = "Meta description"

$metaPara4 = $d.Paragraphs(2)
$metaPara4.Range.InsertAfter(": Read our review of Desperate Dawgs 2 Gigablox and play this exciting game for free today. Features the Gigablox mechanic for larger symbols and three unique bonus rounds.")

# ---------------------------------------------------------------------------
# 2) Remove the now-duplicated bold "Play Desperate Dawgs..." paragraph that
#    used to sit right before the closing italic paragraph.
# ---------------------------------------------------------------------------
$dupPara = $null
for ($i = 2; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ((TrimCR $para.Range.Text) -eq $titleText) {
        $dupPara = $para
    }
}
$delRange = $d.Range($dupPara.Range.Start, $dupPara.Range.End)
$delRange.Delete()

# ---------------------------------------------------------------------------
# 3) Replace the text of the closing italic paragraph with the new
#    image-generation prompt (keeping its italic run formatting intact).
# ---------------------------------------------------------------------------
$oldClosingText = "Read our review of Desperate Dawgs 2 Gigablox and play this exciting game for free today. Features the Gigablox mechanic for larger symbols and three unique bonus rounds."
$newClosingText = "Create a cartoon-style feature image for Desperate Dawgs 2 Gigablox that prominently features a happy Maya warrior wearing glasses. The image should be bright and colorful, and the warrior should be positioned in the center. The warrior should be standing on a Wild West-themed background with symbols from the game, including train carriages, revolvers, and sheriff badges, incorporated into the design. The overall vibe should be fun, exciting, and adventurous, with a clear nod to the Gigablox mechanic and the game's unique features."

$closingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ((TrimCR $para.Range.Text) -eq $oldClosingText) {
        $closingPara = $para
    }
}
$closingRange = $d.Range($closingPara.Range.Start, $closingPara.Range.End - 1)
$closingRange.Text = $newClosingText

Write-Host "done"
